# Adding a sample for a worksheet autofilter
# Appends two rows to the "Snippets" table documenting the new
# Worksheet.autofilter / AutoFilter.apply snippet (excel-worksheet-auto-filter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Snippets")

# Grow the table by two rows; this automatically extends the table ref,
# the autoFilter ref and the sheet dimension (A1:D184 -> A1:D186).
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Fill in the new data, in the same order the values were originally typed,
# so that newly introduced shared strings come out in the same order:
#   autofilter, addAutoFilter, excel-worksheet-auto-filter, AutoFilter, apply
$ws.Range("B185").Value = "autofilter"
$ws.Range("D185").Value = "addAutoFilter"
$ws.Range("C185").Value = "excel-worksheet-auto-filter"
$ws.Range("A186").Value = "AutoFilter"
$ws.Range("B186").Value = "apply"

$ws.Range("A185").Value = "Worksheet"
$ws.Range("C186").Value = "excel-worksheet-auto-filter"
$ws.Range("D186").Value = "addAutoFilter"

# The previous last row (184) loses the stray direct formatting it had.
$ws.Range("A184:D184").ClearFormats()

# Leave the selection where the author ended up after editing.
$ws.Range("O178").Select()
